# Updated pin assignments & diode part # (backorder)
#
# The D1 diode part used in the BOM (row 12) was on backorder, so it is
# swapped for a replacement Schottky diode that's in stock:
#   Comment/LibRef : MBRA340T3G      -> SSA34HE3_A_I
#   Description    : Diode           -> Schottky Diode
#   Footprint      : DIOM5226X220N   -> DIONM5127X229N
#
# A leading "'" forces the new values to stay plain text (matching how the
# existing BOM entries are stored) instead of being re-interpreted/losing
# their formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 = D1 diode entry: A=Comment, B=Description, C=Designator,
# D=Footprint, E=LibRef, F=Quantity.
$ws.Range("A12").Value = "'SSA34HE3_A_I"
$ws.Range("B12").Value = "'Schottky Diode"
$ws.Range("D12").Value = "'DIONM5127X229N"
$ws.Range("E12").Value = "'SSA34HE3_A_I"

# Restore the saved cursor position to the default A1 cell.
$ws.Range("A1").Select()
